$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new BCB-etch process got split into sub-steps (8A-8D) and two brand
# new mesa-size variants (rows 10 & 11) were inserted into the layer
# table. Reproduce the edit in the same order it was originally made so
# the shared-string table comes out the same way.

# 1) Renumber/relabel the BCB-etch-contact step to "8A. ..."
$ws.Range("F9").Value = "8A. BCB ETCH CONT (DC)"

# 2) Insert two new blank rows for the extra mesa sizes
$ws.Rows.Item(10).Resize(2).Insert()

# 3) Fill in the two new rows
$ws.Range("A10").Value = "bcb_etch_mesa_15-21"
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = "DC"
$ws.Range("D10").Value = 9
$ws.Range("E10").Value = 0

$ws.Range("A11").Value = "bcb_etch_mesa_17-23"
$ws.Range("B11").Value = 10
$ws.Range("C11").Value = "DC"
$ws.Range("D11").Value = 9
$ws.Range("E11").Value = 0

# 4) The rows that used to be 10/11/12 are now 12/13/14 - relabel them
$ws.Range("B12").Value = 11
$ws.Range("D12").Value = 8
$ws.Range("F12").Value = "8D. SiNx ON BCB (DC)"

$ws.Range("B13").Value = 12
$ws.Range("D13").Value = "none"
$ws.Range("F13").Value = "9. PADS  (DC)"

$ws.Range("B14").Value = 13
$ws.Range("D14").Value = "none"
$ws.Range("F14").Value = "10. BONDPADS (DC)"

# 5) Former rows 13/14/15 (frame/chip rows) are now 15/16/17
$ws.Range("B15").Value = 14
$ws.Range("B16").Value = 15
$ws.Range("B17").Value = 16

# 6) Go back and fix the wording of the bottom-contact / contact-etch labels
$ws.Range("F6").Value = "5. BOTTOM CONTACT(DC)"
$ws.Range("F7").Value = "6. CONTACT ETCH(DD)"

# 7) Finally fill in the labels for the two new mesa-size rows
$ws.Range("F10").Value = "8B. BCB MESA 15-21(DC)"
$ws.Range("F11").Value = "8C. BCB MESA 17-23(DC)"

# --- Column widths --------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 22.5
$ws.Columns.Item(6).ColumnWidth = 36.5

# --- Selection / active cell ----------------------------------------------
$ws.Range("F14").Select()

# --- Page setup -------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
